$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target data (row -> A:Class, B:Points, C:Title, D:DueDate)
# Row 1: CMST / Quiz / 1010.0 / Sat May 04 00:00:00 CDT 2019
# Row 2: CMST / Essay / 1234.0 / Fri May 10 00:00:00 CDT 2019
# Row 3: MATH / Chapter 10 questions / 3423.0 / Mon May 06 13:54:27 CDT 2019
# Row 4: ITEC / Java final project / 2545.0 / Wed May 15 00:00:00 CDT 2019
# Row 5: Test / test assignment / 1234.0 / Wed May 08 14:34:07 CDT 2019

$ws.Range("A1").Value = "CMST"
$ws.Range("B1").Value = 1010.0
$ws.Range("C1").Value = "Quiz"
$ws.Range("D1").Value = "Sat May 04 00:00:00 CDT 2019"

$ws.Range("A2").Value = "CMST"
$ws.Range("B2").Value = 1234.0
$ws.Range("C2").Value = "Essay"
$ws.Range("D2").Value = "Fri May 10 00:00:00 CDT 2019"

$ws.Range("A3").Value = "MATH"
$ws.Range("B3").Value = 3423.0
$ws.Range("C3").Value = "Chapter 10 questions"
$ws.Range("D3").Value = "Mon May 06 13:54:27 CDT 2019"

$ws.Range("A4").Value = "ITEC"
$ws.Range("B4").Value = 2545.0
$ws.Range("C4").Value = "Java final project"
$ws.Range("D4").Value = "Wed May 15 00:00:00 CDT 2019"

$ws.Range("A5").Value = "Test"
$ws.Range("B5").Value = 1234.0
$ws.Range("C5").Value = "test assignment"
$ws.Range("D5").Value = "Wed May 08 14:34:07 CDT 2019"
